{"js": "const paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\nlet target = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const p = paragraphs.items[i];\n  if (p.text.indexOf(\"Docente(s) Respons\u00e1vel(eis)\") !== -1) {\n    target = p;\n    break;\n  }\n}\n\nif (!target) {\n  throw new Error(\"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\");\n}\n\nconst inserted = target.insertParagraph(\"1285870 - Marcos Villela Barcza\", Word.InsertLocation.after);\ninserted.style = \"ListBullet\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n$found = $null\nforeach ($p in $d.Paragraphs) {\n    if ($p.Range.Text -like \"*Docente(s) Respons\u00e1vel(eis)*\") {\n        $found = $p\n        break\n    }\n}\n\nif ($found -eq $null) {\n    throw \"Could not find the 'Docente(s) Respons\u00e1vel(eis)' paragraph\"\n}\n\n# Insert a brand-new paragraph right after the \"Docente(s) Respons\u00e1vel(eis)\"\n# heading, then populate its text and apply the ListBullet style.\n$r = $found.Range\n$r.Collapse(0)\n$r.InsertParagraphAfter()\n\n$newRange = $found.Next().Range\n$newRange.Text = \"1285870 - Marcos Villela Barcza\"\n$newRange.Style = \"ListBullet\"\n"}
